# Generate Report for Handback
# -----------------------------------------------------------------------
# This script reproduces, via Excel COM automation, the "handback" report
# generation step: the localization status workbook is updated to reflect
# that the zh-cn and de-de targets have been handed back and are in sync
# with en-US.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# -------------------------------------------------------------------
# 1. Overview sheet: the zh-cn / de-de status columns (E, F) move from
#    "Ready for handoff" to "Handed back: in sync with en-US".
# -------------------------------------------------------------------
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# -------------------------------------------------------------------
# 2. zh-cn sheet: fill in the handback info for rows 2 & 3 -
#    - "Status" (col C) text changes along with the shared string above
#    - "Latest Target File" (col I) now links to a.md
#    - "Latest Handback File" (col J) is the generated .xlf
#    - "Latest Handback DateTime" (col K) gets the handback timestamp
# -------------------------------------------------------------------
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Range("I2").Value = "a.md"
$zhcn.Range("I3").Value = "a.md"

$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-28 14:36:57"
$zhcn.Range("K3").Value = "2016-08-28 14:36:57"

$zhcn.Columns.Item(3).ColumnWidth  = 29.166666666666668
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# Add hyperlinks for the new "Latest Target File" entries (I2, I3); the
# existing A2/A3 (a.md / b.md) hyperlinks are left exactly as they are.
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d4d1bb80ac946272d30321f8783690af0491d389/e2e/a.md", "", "", "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d4d1bb80ac946272d30321f8783690af0491d389/e2e/a.md", "", "", "a.md") | Out-Null

# -------------------------------------------------------------------
# 3. de-de sheet: same shape of update as zh-cn, but with its own
#    handback artifact name/timestamp.
# -------------------------------------------------------------------
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("I2").Value = "a.md"
$dede.Range("I3").Value = "a.md"

$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("K2").Value = "2016-08-28 14:37:08"
$dede.Range("K3").Value = "2016-08-28 14:37:08"

$dede.Columns.Item(3).ColumnWidth  = 29.166666666666668
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

# Same as zh-cn: add hyperlinks for I2/I3 only; leave A2/A3 untouched.
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d4d1bb80ac946272d30321f8783690af0491d389/e2e/a.md", "", "", "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d4d1bb80ac946272d30321f8783690af0491d389/e2e/a.md", "", "", "a.md") | Out-Null
